$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.845.27"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "2.532.10"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'590.92"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").Value = "'173.51"
$ws.Range("E6").Value = "  -1.29%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.523"
$ws.Range("E8").Value = "  -1.34%  "

$ws.Range("D9").Value = "2.527.91"
$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("E10").Value = "  -2.21%  "

$ws.Range("E11").Value = "  +1.77%  "

$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").Value = "'5.02"
$ws.Range("E13").Value = "  -2.92%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.018.01"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'26.34"
$ws.Range("E15").Value = "  -1.68%  "

$ws.Range("D16").Value = "'2.45"
$ws.Range("E16").Value = "  +145.73%  "

$ws.Range("D17").Value = "'0.0000176"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").Value = "67.780.06"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("D19").Value = "2.542.95"
$ws.Range("E19").Value = "  +1.97%  "

$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("E21").Value = "  -2.14%  "

$ws.Range("D22").Value = "'369.98"
$ws.Range("E22").Value = "  +2.96%  "

$ws.Range("E23").Value = "  -1.45%  "

$ws.Range("D24").Value = "'4.56"
$ws.Range("E24").Value = "  -2.10%  "

$ws.Range("D25").Value = "'71.84"
$ws.Range("E25").Value = "  +2.45%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'1.91"
$ws.Range("E27").Value = "  -4.38%  "

$ws.Range("D28").Value = "'9.96"
$ws.Range("E28").Value = "  -2.61%  "

$ws.Range("D30").Value = "0.0₃0966"
$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("D31").Value = "'538.73"
$ws.Range("E31").Value = "  -2.37%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  -2.88%  "

$ws.Range("D34").Value = "'1.86"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").Value = "'160.09"
$ws.Range("E37").Value = "  +1.29%  "

$ws.Range("E38").Value = "  -2.94%  "

$ws.Range("E39").Value = "  +2.03%  "

$ws.Range("D40").Value = "'18.61"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.350"
$ws.Range("E41").Value = "  -1.80%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.14"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("D43").Value = "'1.77"
$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("E44").Value = "  -2.15%  "

$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "'39.34"
$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").Value = "0.0₆0283"
$ws.Range("E47").Value = "  +1.65%  "

$ws.Range("D48").Value = "'147.85"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("D49").Value = "'3.70"
$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").Value = "'0.549"
$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("E51").Value = "  +0.38%  "
